# Apply "dSF" (column F) corrections for the martínez_nick 2024 save-data sheet.
# Column F ("dSF") values are being repulled/recalculated; only column F
# changes -- column E ("dS0") and all other columns stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (6)
$updates = @{
    2  = -1
    3  = 8
    4  = 9
    5  = 1
    6  = 11
    8  = -2
    9  = 1
    10 = -3
    12 = 7
    13 = -3
    14 = 5
    18 = 1
    19 = -2
    20 = 4
    21 = 2
    26 = -1
    27 = 3
    28 = -2
    30 = -1
    32 = 4
    33 = -2
    34 = 4
    35 = -2
    36 = -3
    38 = -5
    39 = -5
    40 = 2
    42 = -3
    44 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
